$p = $ppt.ActivePresentation

# Find the title placeholder that contains the "Path ... (Postman)" line
# (slide 12 in this deck) instead of hard-coding indices.
$tr = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text.Contains("наименование запроса")) {
                $tr = $shp.TextFrame.TextRange
            }
        }
    }
}

# Locate the "Path ... (Postman)" line. We rebuild the fragment that used to
# read:  Path  - наименование запроса (Postman)
# into:  Path  - endpoint (Postman)
# as four runs:  " " / "- " / "endpoint" (sz 2400) / " "  followed by the
# untouched "(Postman)" run.

$full    = $tr.Text
$pathIdx = $full.IndexOf("Path")        # 0-based index of the "P" in "Path"
$dash1   = $pathIdx + 1 + 5             # 1-based index of the char right after "Path "

# Old layout starting at $dash1 (24 chars total): " - наименование запроса "
#   run A (3 chars):  " - "
#   run B (13 chars): "наименование "
#   run C (8 chars):  "запроса "

# 1) Shrink run A from " - " down to just " " (keeps its rPr: lang en-US, dirty 0).
$runA = $tr.Characters($dash1, 3)
$runA.Text = " "

# 2) Insert the new "- " run right after that single space. It inherits the
#    neighbouring lang="en-US" formatting.
$afterSpace = $tr.Characters($dash1, 1)
$afterSpace.InsertAfter("- ") | Out-Null

# 3) Insert the new "endpoint" run after "- ", still lang="en-US", then resize
#    it down to 24pt to match the target.
$dashRun = $tr.Characters($dash1 + 1, 2)
$dashRun.InsertAfter("endpoint") | Out-Null
$endpointRun = $tr.Characters($dash1 + 3, 8)
$endpointRun.Font.Size = 24

# 4) The old "наименование " (run B) and "запроса " (run C) runs are still
#    sitting right after "endpoint". Drop run B entirely and shrink run C
#    down to a single space so the space keeps run C's own rPr
#    (lang="ru-RU" smtClean="0").
$runBStart = $dash1 + 3 + 8
$runB = $tr.Characters($runBStart, 13)
$runB.Delete()

$runC = $tr.Characters($runBStart, 8)
$runC.Text = " "
